$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "57.828.72"
$ws.Cells.Item(2, 5).Value = "  +2.59%  "
$ws.Cells.Item(3, 4).Value = "2.351.80"
$ws.Cells.Item(3, 5).Value = "  +1.67%  "
$ws.Cells.Item(4, 5).Value = "  -0.19%  "
$ws.Cells.Item(5, 4).Value = "'547.14"
$ws.Cells.Item(5, 5).Value = "  +6.25%  "
$ws.Cells.Item(6, 4).Value = "'134.99"
$ws.Cells.Item(6, 5).Value = "  +2.76%  "
$ws.Cells.Item(7, 5).Value = "  -0.03%  "
$ws.Cells.Item(8, 4).Value = "'0.572"
$ws.Cells.Item(8, 5).Value = "  +7.82%  "
$ws.Cells.Item(9, 4).Value = "2.349.57"
$ws.Cells.Item(9, 5).Value = "  +1.43%  "
$ws.Cells.Item(10, 5).Value = "  +1.92%  "
$ws.Cells.Item(11, 5).Value = "  +3.70%  "
$ws.Cells.Item(12, 5).Value = "  +0.04%  "
$ws.Cells.Item(13, 5).Value = "  +6.83%  "
$ws.Cells.Item(14, 4).Value = "2.763.97"
$ws.Cells.Item(14, 5).Value = "  +1.12%  "
$ws.Cells.Item(15, 5).Value = "  +0.77%  "
$ws.Cells.Item(16, 4).Value = "57.788.32"
$ws.Cells.Item(16, 5).Value = "  +2.46%  "
$ws.Cells.Item(17, 5).Value = "  +1.12%  "
$ws.Cells.Item(18, 4).Value = "2.356.51"
$ws.Cells.Item(18, 5).Value = "  +0.98%  "
$ws.Cells.Item(19, 4).Value = "'10.65"
$ws.Cells.Item(19, 5).Value = "  +2.70%  "
$ws.Cells.Item(20, 4).Value = "'334.93"
$ws.Cells.Item(20, 5).Value = "  +1.75%  "
$ws.Cells.Item(21, 4).Value = "'4.24"
$ws.Cells.Item(21, 5).Value = "  +2.38%  "
$ws.Cells.Item(22, 4).Value = "'6.72"
$ws.Cells.Item(22, 5).Value = "  +0.52%  "
$ws.Cells.Item(23, 5).Value = "  +0.16%  "
$ws.Cells.Item(24, 4).Value = "'5.59"
$ws.Cells.Item(24, 5).Value = "  +0.94%  "
$ws.Cells.Item(25, 4).Value = "'62.17"
$ws.Cells.Item(25, 5).Value = "  +2.07%  "
$ws.Cells.Item(26, 5).Value = "  +2.07%  "
$ws.Cells.Item(27, 4).Value = "'8.52"
$ws.Cells.Item(27, 5).Value = "  -1.23%  "
$ws.Cells.Item(28, 4).Value = "'0.999"
$ws.Cells.Item(28, 5).Value = "  +0.42%  "
$ws.Cells.Item(29, 5).Value = "  +6.46%  "
$ws.Cells.Item(30, 4).Value = "'1.78"
$ws.Cells.Item(30, 5).Value = "  +5.73%  "
$ws.Cells.Item(31, 4).Value = "'170.27"
$ws.Cells.Item(31, 5).Value = "  +1.17%  "
$ws.Cells.Item(32, 4).Value = "0.0₃0735"
$ws.Cells.Item(32, 5).Value = "  +2.52%  "
$ws.Cells.Item(33, 4).Value = "'6.16"
$ws.Cells.Item(33, 5).Value = "  +0.54%  "
$ws.Cells.Item(34, 4).Value = "'1.04"
$ws.Cells.Item(34, 5).Value = "  +17.77%  "
$ws.Cells.Item(35, 4).Value = "'18.52"
$ws.Cells.Item(35, 5).Value = "  +1.39%  "
$ws.Cells.Item(36, 5).Value = "  +0.03%  "
$ws.Cells.Item(37, 4).Value = "'0.999"
$ws.Cells.Item(37, 5).Value = "  +0.01%  "
$ws.Cells.Item(38, 4).Value = "'4.18"
$ws.Cells.Item(38, 5).Value = "  +6.63%  "
$ws.Cells.Item(39, 5).Value = "  +1.01%  "
$ws.Cells.Item(40, 5).Value = "  +2.99%  "
$ws.Cells.Item(41, 4).Value = "'39.12"
$ws.Cells.Item(41, 5).Value = "  +1.40%  "
$ws.Cells.Item(42, 4).Value = "'148.06"
$ws.Cells.Item(42, 5).Value = "  -0.29%  "
$ws.Cells.Item(43, 5).Value = "  +0.98%  "
$ws.Cells.Item(44, 4).Value = "'3.63"
$ws.Cells.Item(44, 5).Value = "  +1.48%  "
$ws.Cells.Item(45, 4).Value = "'285.82"
$ws.Cells.Item(45, 5).Value = "  +0.64%  "
$ws.Cells.Item(46, 4).Value = "'0.0944"
$ws.Cells.Item(46, 5).Value = "  +2.02%  "
$ws.Cells.Item(47, 4).Value = "'19.22"
$ws.Cells.Item(47, 5).Value = "  +6.50%  "
$ws.Cells.Item(48, 5).Value = "  +2.52%  "
$ws.Cells.Item(49, 4).Value = "'0.562"
$ws.Cells.Item(49, 5).Value = "  +1.64%  "
$ws.Cells.Item(50, 5).Value = "  +1.67%  "
$ws.Cells.Item(51, 2).Value = "Polygon"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(51, 4).Value = "'0.385"
$ws.Cells.Item(51, 5).Value = "  +7.40%  "

# Cells whose numeric-looking text needed a quote-prefix to stay text;
# reset their style back to Normal so no stray quotePrefix style lingers.
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(51, 4).Style = "Normal"
